# Actualización automática desde tarea programada
# Updates the timestamp precision on row 3 and appends a new sensor
# reading as row 4 (mirrors the style/format used by the existing data
# rows 2 and 3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: refresh the serial-date value (tiny precision correction) ---
$ws.Cells.Item(3, 1).Value = 45873.37521892361

# --- Row 4: new data row appended by the scheduled task ---
$ws.Cells.Item(4, 1).Value = 45873.41691728914
$ws.Cells.Item(4, 1).NumberFormat = $ws.Cells.Item(3, 1).NumberFormat

$ws.Cells.Item(4, 2).Value = 2025
$ws.Cells.Item(4, 3).Value = 15
$ws.Cells.Item(4, 4).Value = 17
$ws.Cells.Item(4, 5).Value = 83.41
$ws.Cells.Item(4, 6).Value = 472.85
$ws.Cells.Item(4, 7).Value = 7.04
$ws.Cells.Item(4, 8).Value = "ESE"
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 10).Value = "10:00:21"
